$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44161
$ws.Range("K2").Value = 'Bing'
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 160
$ws.Range("N2").Value = 39000
$ws.Range("O2").Value = 40000
$ws.Range("P2").Value = 39500
$ws.Range("Q2").Value = '$/caja 20 kilos'
$ws.Range("R2").Value = 'Provincia de Curicó'
$ws.Range("S2").Value = 1975
$ws.Range("T2").Value = 20

# Row 3
$ws.Range("D3").Value = 44208
$ws.Range("K3").Value = 'Lapins'
$ws.Range("L3").Value = 'Segunda'
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 10500
$ws.Range("O3").Value = 11000
$ws.Range("P3").Value = 10750
$ws.Range("Q3").Value = '$/bandeja 12 kilos'
$ws.Range("R3").Value = 'Provincia de Curicó'
$ws.Range("S3").Value = 896
$ws.Range("T3").Value = 12

# Row 4
$ws.Range("D4").Value = 44537
$ws.Range("K4").Value = 'Brooks'
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 29000
$ws.Range("O4").Value = 30000
$ws.Range("P4").Value = 29500
$ws.Range("Q4").Value = '$/caja 20 kilos'
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 1475
$ws.Range("T4").Value = 20

# Row 5
$ws.Range("D5").Value = 44922
$ws.Range("K5").Value = 'Bing'
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 5000
$ws.Range("O5").Value = 6000
$ws.Range("P5").Value = 5500
$ws.Range("Q5").Value = '$/bandeja 10 kilos'
$ws.Range("R5").Value = 'Región del Maule'
$ws.Range("S5").Value = 550
$ws.Range("T5").Value = 10

# Row 6
$ws.Range("D6").Value = 44532
$ws.Range("K6").Value = 'Brooks'
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 400
$ws.Range("N6").Value = 27000
$ws.Range("O6").Value = 28000
$ws.Range("P6").Value = 27500
$ws.Range("Q6").Value = '$/bandeja 12 kilos'
$ws.Range("R6").Value = 'Región de O''Higgins'
$ws.Range("S6").Value = 2292
$ws.Range("T6").Value = 12

# Row 7
$ws.Range("D7").Value = 44210
$ws.Range("K7").Value = 'Rainier'
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 250
$ws.Range("N7").Value = 21000
$ws.Range("O7").Value = 22000
$ws.Range("P7").Value = 21500
$ws.Range("Q7").Value = '$/caja 18 kilos'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 1194
$ws.Range("T7").Value = 18

# Row 8
$ws.Range("D8").Value = 44914
$ws.Range("K8").Value = 'Brooks'
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 700
$ws.Range("N8").Value = 7000
$ws.Range("O8").Value = 8000
$ws.Range("P8").Value = 7429
$ws.Range("Q8").Value = '$/bandeja 10 kilos'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 743
$ws.Range("T8").Value = 10

# Row 9
$ws.Range("D9").Value = 44914
$ws.Range("K9").Value = 'Lapins'
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 550
$ws.Range("N9").Value = 7000
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 7455
$ws.Range("Q9").Value = '$/bandeja 10 kilos'
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 746
$ws.Range("T9").Value = 10

# Row 10
$ws.Range("D10").Value = 44571
$ws.Range("K10").Value = 'Brooks'
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 400
$ws.Range("N10").Value = 8500
$ws.Range("O10").Value = 9000
$ws.Range("P10").Value = 8750
$ws.Range("Q10").Value = '$/bandeja 10 kilos'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 875
$ws.Range("T10").Value = 10

# Row 11
$ws.Range("D11").Value = 44229
$ws.Range("K11").Value = 'Santina'
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 250
$ws.Range("N11").Value = 6500
$ws.Range("O11").Value = 7000
$ws.Range("P11").Value = 6750
$ws.Range("Q11").Value = '$/bandeja 5 kilos'
$ws.Range("R11").Value = 'Provincia de Curicó'
$ws.Range("S11").Value = 1350
$ws.Range("T11").Value = 5

# Row 12
$ws.Range("D12").Value = 44943
$ws.Range("K12").Value = 'Santina'
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 600
$ws.Range("N12").Value = 14000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 14333
$ws.Range("Q12").Value = '$/caja 15 kilos'
$ws.Range("R12").Value = 'Región del Maule'
$ws.Range("S12").Value = 956
$ws.Range("T12").Value = 15

# Row 13
$ws.Range("D13").Value = 44921
$ws.Range("K13").Value = 'Bing'
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 320
$ws.Range("N13").Value = 7500
$ws.Range("O13").Value = 8000
$ws.Range("P13").Value = 7781
$ws.Range("Q13").Value = '$/bandeja 10 kilos'
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("S13").Value = 778
$ws.Range("T13").Value = 10

# Row 14
$ws.Range("D14").Value = 44557
$ws.Range("K14").Value = 'Lapins'
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 250
$ws.Range("N14").Value = 9000
$ws.Range("O14").Value = 10000
$ws.Range("P14").Value = 9500
$ws.Range("Q14").Value = '$/bandeja 10 kilos'
$ws.Range("R14").Value = 'Provincia de Curicó'
$ws.Range("S14").Value = 950
$ws.Range("T14").Value = 10

# Row 15
$ws.Range("D15").Value = 44931
$ws.Range("K15").Value = 'Lapins'
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 250
$ws.Range("N15").Value = 6000
$ws.Range("O15").Value = 6500
$ws.Range("P15").Value = 6250
$ws.Range("Q15").Value = '$/bandeja 10 kilos'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 625
$ws.Range("T15").Value = 10

# Row 16
$ws.Range("D16").Value = 44931
$ws.Range("K16").Value = 'Lapins'
$ws.Range("L16").Value = 'Segunda'
$ws.Range("M16").Value = 400
$ws.Range("N16").Value = 3000
$ws.Range("O16").Value = 3300
$ws.Range("P16").Value = 3150
$ws.Range("Q16").Value = '$/bandeja 5 kilos'
$ws.Range("R16").Value = 'Región de O''Higgins'
$ws.Range("S16").Value = 630
$ws.Range("T16").Value = 5

# Row 17
$ws.Range("D17").Value = 44908
$ws.Range("K17").Value = 'Rainier'
$ws.Range("L17").Value = 'Segunda'
$ws.Range("M17").Value = 250
$ws.Range("N17").Value = 15000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 15600
$ws.Range("Q17").Value = '$/caja 10 kilos'
$ws.Range("R17").Value = 'Región de O''Higgins'
$ws.Range("S17").Value = 1560
$ws.Range("T17").Value = 10

# Row 18
$ws.Range("D18").Value = 44175
$ws.Range("K18").Value = 'Rainier'
$ws.Range("L18").Value = 'Segunda'
$ws.Range("M18").Value = 270
$ws.Range("N18").Value = 25000
$ws.Range("O18").Value = 26000
$ws.Range("P18").Value = 25500
$ws.Range("Q18").Value = '$/caja 18 kilos'
$ws.Range("R18").Value = 'Región de O''Higgins'
$ws.Range("S18").Value = 1417
$ws.Range("T18").Value = 18

# Row 19
$ws.Range("D19").Value = 44580
$ws.Range("K19").Value = 'Sweet Heart'
$ws.Range("L19").Value = 'Segunda'
$ws.Range("M19").Value = 300
$ws.Range("N19").Value = 7000
$ws.Range("O19").Value = 8000
$ws.Range("P19").Value = 7500
$ws.Range("Q19").Value = '$/bandeja 10 kilos'
$ws.Range("R19").Value = 'Región de O''Higgins'
$ws.Range("S19").Value = 750
$ws.Range("T19").Value = 10

# Row 20
$ws.Range("D20").Value = 44901
$ws.Range("K20").Value = 'Bing'
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 500
$ws.Range("N20").Value = 12000
$ws.Range("O20").Value = 13000
$ws.Range("P20").Value = 12500
$ws.Range("Q20").Value = '$/caja 15 kilos'
$ws.Range("R20").Value = 'Región de O''Higgins'
$ws.Range("S20").Value = 833
$ws.Range("T20").Value = 15

# Row 21
$ws.Range("D21").Value = 44901
$ws.Range("K21").Value = 'Lapins'
$ws.Range("L21").Value = 'Primera'
$ws.Range("M21").Value = 500
$ws.Range("N21").Value = 12000
$ws.Range("O21").Value = 13000
$ws.Range("P21").Value = 12500
$ws.Range("Q21").Value = '$/caja 15 kilos'
$ws.Range("R21").Value = 'Región de O''Higgins'
$ws.Range("S21").Value = 833
$ws.Range("T21").Value = 15

# Row 22
$ws.Range("D22").Value = 44568
$ws.Range("K22").Value = 'Santina'
$ws.Range("L22").Value = 'Segunda'
$ws.Range("M22").Value = 200
$ws.Range("N22").Value = 15000
$ws.Range("O22").Value = 16000
$ws.Range("P22").Value = 15500
$ws.Range("Q22").Value = '$/bandeja 12 kilos'
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("S22").Value = 1292
$ws.Range("T22").Value = 12

# Row 23
$ws.Range("D23").Value = 44917
$ws.Range("K23").Value = 'Bing'
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 400
$ws.Range("N23").Value = 5000
$ws.Range("O23").Value = 6000
$ws.Range("P23").Value = 5625
$ws.Range("Q23").Value = '$/bandeja 10 kilos'
$ws.Range("R23").Value = 'Región de O''Higgins'
$ws.Range("S23").Value = 562
$ws.Range("T23").Value = 10

# Row 24
$ws.Range("D24").Value = 44917
$ws.Range("K24").Value = 'Santina'
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 400
$ws.Range("N24").Value = 5000
$ws.Range("O24").Value = 6000
$ws.Range("P24").Value = 5500
$ws.Range("Q24").Value = '$/bandeja 10 kilos'
$ws.Range("R24").Value = 'Región de O''Higgins'
$ws.Range("S24").Value = 550
$ws.Range("T24").Value = 10
